$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp label in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 13:56"

# Row 4
$ws.Range("B4").Value = 2330908
$ws.Range("C4").Value = 330
$ws.Range("E4").Value = 1235854
$ws.Range("G4").Value = 19
$ws.Range("H4").Value = 121999

# Row 7
$ws.Range("B7").Value = 412955
$ws.Range("C7").Value = 1228
$ws.Range("D7").Value = 228605
$ws.Range("E7").Value = 171058
$ws.Range("G7").Value = 15
$ws.Range("H7").Value = 13292

# Row 13
$ws.Range("B13").Value = 204952
$ws.Range("C13").Value = 2368
$ws.Range("D13").Value = 163591
$ws.Range("E13").Value = 31738
$ws.Range("G13").Value = 116
$ws.Range("H13").Value = 9623

# Row 23
$ws.Range("B23").Value = 87369
$ws.Range("C23").Value = 881
$ws.Range("D23").Value = 68319
$ws.Range("E23").Value = 18952
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 98

# Row 27
$ws.Range("B27").Value = 58505
$ws.Range("C27").Value = 569
$ws.Range("D27").Value = 37666
$ws.Range("E27").Value = 20493
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 346

# Row 33
$ws.Range("B33").Value = 44925
$ws.Range("C33").Value = 392
$ws.Range("D33").Value = 32415
$ws.Range("E33").Value = 12208
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 302

# Row 36
$ws.Range("B36").Value = 39650
$ws.Range("C36").Value = 505
$ws.Range("D36").Value = 31240
$ws.Range("E36").Value = 8084
$ws.Range("G36").Value = 7
$ws.Range("H36").Value = 326

# Row 55
$ws.Range("B55").Value = 17341
$ws.Range("C55").Value = 18
$ws.Range("D55").Value = 16197
$ws.Range("E55").Value = 454
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 690

# Row 69
$ws.Range("A69").Value = "Nepal"
$ws.Range("B69").Value = 9026
$ws.Range("C69").Value = 421
$ws.Range("D69").Value = 1772
$ws.Range("E69").Value = 7231
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 23

# Row 70
$ws.Range("A70").Value = "Noruega"
$ws.Range("B70").Value = 8742
$ws.Range("D70").Value = 8138
$ws.Range("E70").Value = 360
$ws.Range("H70").Value = 244

# Row 77
$ws.Range("A77").Value = "Senegal"
$ws.Range("B77").Value = 5888
$ws.Range("C77").Value = 105
$ws.Range("D77").Value = 3919
$ws.Range("E77").Value = 1885
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 84

# Row 78
$ws.Range("A78").Value = "Consejo Danes para los Refugiados"
$ws.Range("B78").Value = 5826
$ws.Range("C78").Value = 154
$ws.Range("D78").Value = 841
$ws.Range("E78").Value = 4855
$ws.Range("G78").Value = 5
$ws.Range("H78").Value = 130

# Row 118
$ws.Range("B118").Value = 1520
$ws.Range("C118").Value = 1
$ws.Range("D118").Value = 1376
$ws.Range("E118").Value = 35

# Row 123
$ws.Range("B123").Value = 1157
$ws.Range("C123").Value = 1
$ws.Range("D123").Value = 1020
$ws.Range("E123").Value = 87

# Row 131
$ws.Range("D131").Value = 814
$ws.Range("E131").Value = 34

# Row 139
$ws.Range("D139").Value = 258
$ws.Range("E139").Value = 461

# Row 145
$ws.Range("B145").Value = 665
$ws.Range("C145").Value = 1
$ws.Range("E145").Value = 40

# Row 202
$ws.Range("A202").Value = "Dominica"

# Row 203
$ws.Range("A203").Value = "Fiyi"

# Row 208
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

# Row 209
$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0

# Row 213
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

# Row 214
$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
